$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Etar vs Spartak Pleven) - updated odds
$ws.Cells.Item(7, 7).Value = 1.4
$ws.Cells.Item(7, 8).Value = 4
$ws.Cells.Item(7, 9).Value = 6.7
$ws.Cells.Item(7, 10).Value = 1.78
$ws.Cells.Item(7, 11).Value = 1.83
$ws.Cells.Item(7, 14).Value = 5.3
$ws.Cells.Item(7, 15).Value = 5.4
$ws.Cells.Item(7, 16).Value = 7
$ws.Cells.Item(7, 17).Value = 7.6
$ws.Cells.Item(7, 18).Value = 10
$ws.Cells.Item(7, 19).Value = 23
$ws.Cells.Item(7, 21).Value = 7
$ws.Cells.Item(7, 22).Value = 16.5
$ws.Cells.Item(7, 23).Value = 75
$ws.Cells.Item(7, 25).Value = 13.5
$ws.Cells.Item(7, 26).Value = 32
$ws.Cells.Item(7, 27).Value = 17.5
$ws.Cells.Item(7, 28).Value = 110
$ws.Cells.Item(7, 29).Value = 60
$ws.Cells.Item(7, 30).Value = 55

# Row 10 - updated odds
$ws.Cells.Item(10, 7).Value = 1.48
$ws.Cells.Item(10, 8).Value = 4.75
$ws.Cells.Item(10, 10).Value = 1.6
$ws.Cells.Item(10, 11).Value = 2.3
$ws.Cells.Item(10, 19).Value = 21
$ws.Cells.Item(10, 20).Value = 17
$ws.Cells.Item(10, 33).Value = 1.18
$ws.Cells.Item(10, 34).Value = 4.5

# Row 12 - updated odds
$ws.Cells.Item(12, 7).Value = 2.5
$ws.Cells.Item(12, 9).Value = 2.77
$ws.Cells.Item(12, 10).Value = 2.37
$ws.Cells.Item(12, 11).Value = 1.45
$ws.Cells.Item(12, 14).Value = 6.3
$ws.Cells.Item(12, 15).Value = 10.75
$ws.Cells.Item(12, 16).Value = 10.5
$ws.Cells.Item(12, 17).Value = 27
$ws.Cells.Item(12, 18).Value = 26
$ws.Cells.Item(12, 20).Value = 6.7
$ws.Cells.Item(12, 21).Value = 6
$ws.Cells.Item(12, 25).Value = 6.7
$ws.Cells.Item(12, 26).Value = 12.5
$ws.Cells.Item(12, 27).Value = 11
$ws.Cells.Item(12, 28).Value = 32
$ws.Cells.Item(12, 29).Value = 29
$ws.Cells.Item(12, 34).Value = 2.32
$ws.Cells.Item(12, 35).Value = 2.05

# Row 14 - updated odds
$ws.Cells.Item(14, 36).Value = 1.67

# Row 15 - updated odds
$ws.Cells.Item(15, 7).Value = 1.31
$ws.Cells.Item(15, 8).Value = 4.1
$ws.Cells.Item(15, 9).Value = 9.5
$ws.Cells.Item(15, 10).Value = 1.87
$ws.Cells.Item(15, 11).Value = 1.75
$ws.Cells.Item(15, 12).Value = 1.39
$ws.Cells.Item(15, 13).Value = 2.42
$ws.Cells.Item(15, 14).Value = 4.65
$ws.Cells.Item(15, 15).Value = 4.6
$ws.Cells.Item(15, 16).Value = 7.4
$ws.Cells.Item(15, 17).Value = 6.4
$ws.Cells.Item(15, 19).Value = 29
$ws.Cells.Item(15, 20).Value = 9
$ws.Cells.Item(15, 21).Value = 7.4
$ws.Cells.Item(15, 22).Value = 20
$ws.Cells.Item(15, 23).Value = 110
$ws.Cells.Item(15, 25).Value = 16.5
$ws.Cells.Item(15, 26).Value = 50
$ws.Cells.Item(15, 27).Value = 24
$ws.Cells.Item(15, 28).Value = 200
$ws.Cells.Item(15, 29).Value = 100
$ws.Cells.Item(15, 30).Value = 90

# Row 17 - updated odds
$ws.Cells.Item(17, 7).Value = 2.8
$ws.Cells.Item(17, 8).Value = 3.25
$ws.Cells.Item(17, 9).Value = 2.55
$ws.Cells.Item(17, 10).Value = 2.2
$ws.Cells.Item(17, 11).Value = 1.65
$ws.Cells.Item(17, 25).Value = 7.5
$ws.Cells.Item(17, 29).Value = 21
$ws.Cells.Item(17, 35).Value = 1.91
$ws.Cells.Item(17, 36).Value = 1.91

# Row 18 - updated odds
$ws.Cells.Item(18, 35).Value = 1.7

# Row 19 - updated odds
$ws.Cells.Item(19, 7).Value = 1.37
$ws.Cells.Item(19, 10).Value = 2.08
$ws.Cells.Item(19, 11).Value = 1.73

# Row 20 - updated odds
$ws.Cells.Item(20, 7).Value = 1.39
$ws.Cells.Item(20, 8).Value = 3.75
$ws.Cells.Item(20, 9).Value = 7.5
$ws.Cells.Item(20, 10).Value = 2.15
$ws.Cells.Item(20, 11).Value = 1.67
$ws.Cells.Item(20, 14).Value = 5.5
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(20, 20).Value = 8
$ws.Cells.Item(20, 26).Value = 41
$ws.Cells.Item(20, 27).Value = 26
$ws.Cells.Item(20, 28).Value = 101
$ws.Cells.Item(20, 31).Value = 1.07
$ws.Cells.Item(20, 32).Value = 9
$ws.Cells.Item(20, 33).Value = 1.36
$ws.Cells.Item(20, 34).Value = 3
$ws.Cells.Item(20, 35).Value = 2.38
$ws.Cells.Item(20, 36).Value = 1.53

# Row 21 - updated odds
$ws.Cells.Item(21, 7).Value = 1.42
$ws.Cells.Item(21, 10).Value = 2.2
$ws.Cells.Item(21, 11).Value = 1.65

# Row 22 - updated odds
$ws.Cells.Item(22, 9).Value = 1.5
$ws.Cells.Item(22, 10).Value = 1.93
$ws.Cells.Item(22, 11).Value = 1.88
$ws.Cells.Item(22, 31).Value = 1.05
$ws.Cells.Item(22, 32).Value = 11

# Row 23 - updated odds
$ws.Cells.Item(23, 11).Value = 1.48
$ws.Cells.Item(23, 14).Value = 6.5
$ws.Cells.Item(23, 15).Value = 11
$ws.Cells.Item(23, 18).Value = 26
$ws.Cells.Item(23, 20).Value = 6.5
$ws.Cells.Item(23, 23).Value = 67
$ws.Cells.Item(23, 27).Value = 12
$ws.Cells.Item(23, 29).Value = 29

# Row 25 - updated odds
$ws.Cells.Item(25, 10).Value = 2.6
$ws.Cells.Item(25, 11).Value = 1.48

# Row 26 (San Lorenzo vs Guarani de Fram) - newly populated odds
$ws.Cells.Item(26, 7).Value = 1.98
$ws.Cells.Item(26, 8).Value = 3.05
$ws.Cells.Item(26, 9).Value = 3.95
$ws.Cells.Item(26, 10).Value = 2.27
$ws.Cells.Item(26, 11).Value = 1.5
$ws.Cells.Item(26, 12).Value = 1.53
$ws.Cells.Item(26, 13).Value = 2.2
$ws.Cells.Item(26, 14).Value = 5.6
$ws.Cells.Item(26, 15).Value = 8.25
$ws.Cells.Item(26, 16).Value = 9
$ws.Cells.Item(26, 17).Value = 17.5
$ws.Cells.Item(26, 18).Value = 18.5
$ws.Cells.Item(26, 19).Value = 37
$ws.Cells.Item(26, 20).Value = 6.9
$ws.Cells.Item(26, 21).Value = 6
$ws.Cells.Item(26, 22).Value = 18
$ws.Cells.Item(26, 23).Value = 110
$ws.Cells.Item(26, 25).Value = 8.75
$ws.Cells.Item(26, 26).Value = 20
$ws.Cells.Item(26, 27).Value = 14
$ws.Cells.Item(26, 28).Value = 65
$ws.Cells.Item(26, 29).Value = 45
$ws.Cells.Item(26, 30).Value = 60
$ws.Cells.Item(26, 33).Value = 1.44
$ws.Cells.Item(26, 34).Value = 2.4
$ws.Cells.Item(26, 35).Value = 2
$ws.Cells.Item(26, 36).Value = 1.65

# Row 27 (Guairena vs Sp. Carapegua) - newly populated odds
$ws.Cells.Item(27, 7).Value = 2.37
$ws.Cells.Item(27, 8).Value = 2.92
$ws.Cells.Item(27, 9).Value = 3.05
$ws.Cells.Item(27, 10).Value = 2.32
$ws.Cells.Item(27, 11).Value = 1.47
$ws.Cells.Item(27, 12).Value = 1.55
$ws.Cells.Item(27, 13).Value = 2.15
$ws.Cells.Item(27, 14).Value = 6.1
$ws.Cells.Item(27, 15).Value = 10.25
$ws.Cells.Item(27, 16).Value = 9.75
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = 24
$ws.Cells.Item(27, 19).Value = 40
$ws.Cells.Item(27, 20).Value = 6.6
$ws.Cells.Item(27, 21).Value = 5.8
$ws.Cells.Item(27, 22).Value = 17
$ws.Cells.Item(27, 23).Value = 100
$ws.Cells.Item(27, 25).Value = 7.5
$ws.Cells.Item(27, 26).Value = 14.5
$ws.Cells.Item(27, 27).Value = 11.25
$ws.Cells.Item(27, 28).Value = 40
$ws.Cells.Item(27, 29).Value = 32
$ws.Cells.Item(27, 30).Value = 45
$ws.Cells.Item(27, 33).Value = 1.45
$ws.Cells.Item(27, 34).Value = 2.37
$ws.Cells.Item(27, 35).Value = 1.98
$ws.Cells.Item(27, 36).Value = 1.65
